$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": a new day ("07-nov") is inserted as a column just before
# the old "DJ" column (which held "01-oct."), shifting DJ:EN -> DK:EO.
# The freshly inserted column gets the new date header in row 1 and "-"
# placeholders (no data published yet) in the 24 data rows below it.
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

$wsPrix.Columns("DJ").Insert()

$wsPrix.Range("DJ1").Value = "07-nov"

for ($row = 2; $row -le 25; $row++) {
    $wsPrix.Range("DJ$row").Value = "-"
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append the new day's closing price as row 143.
# The date must stay literal text ("2025-11-05"), like the rows above it, so
# it is entered with a leading apostrophe to stop Excel from reinterpreting
# it as a date serial; the style is then reset to "Normal" so the apostrophe
# doesn't leave a stray quote-prefixed cell style behind.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A143").Value = "'2025-11-05"
$wsGaz.Range("A143").Style = "Normal"
$wsGaz.Range("B143").Value = 30.425

# ---------------------------------------------------------------------------
# Sheet "CO2": same update, append row 143.
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A143").Value = "'2025-11-05"
$wsCo2.Range("A143").Style = "Normal"
$wsCo2.Range("B143").Value = 81.18000000000001
